$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.659.30"
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = "'2.326.09"
$ws.Range('E3').Value = '  +4.23%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'271.79"
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').Value = "'95.73"
$ws.Range('E6').Value = '  +8.87%  '
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = "'0.620"
$ws.Range('E9').Value = '  +3.21%  '
$ws.Range('D10').Value = "'44.86"
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').Value = "'0.0946"
$ws.Range('E11').Value = '  +3.00%  '
$ws.Range('D12').Value = "'8.07"
$ws.Range('E12').Value = '  +6.62%  '
$ws.Range('D13').Value = "'0.105"
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').Value = "'2.675.31"
$ws.Range('E14').Value = '  +4.31%  '
$ws.Range('D15').Value = "'15.63"
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('E16').Value = '  +8.06%  '
$ws.Range('D17').Value = "'2.329.02"
$ws.Range('E17').Value = '  +3.59%  '
$ws.Range('D18').Value = "'43.630.08"
$ws.Range('E19').Value = '  +3.83%  '
$ws.Range('D20').Value = "'6.32"
$ws.Range('E20').Value = '  +6.03%  '
$ws.Range('D21').Value = "'71.96"
$ws.Range('E21').Value = '  +2.53%  '
$ws.Range('D22').Value = "'238.03"
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('E23').Value = '  -3.76%  '
$ws.Range('D24').Value = "'9.52"
$ws.Range('E24').Value = '  +9.64%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = "'2.55"
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('E27').Value = '  +5.06%  '
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = "'38.67"
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('D31').Value = "'22.47"
$ws.Range('E31').Value = '  +8.35%  '
$ws.Range('D32').Value = "'172.74"
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('D34').Value = "'5.48"
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('E35').Value = '  +3.69%  '
$ws.Range('D36').Value = "'0.0358"
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('E37').Value = '  -4.00%  '
$ws.Range('D38').Value = "'4.37"
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').Value = "'3.40"
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('E40').Value = '  +9.41%  '
$ws.Range('E41').Value = '  +11.62%  '
$ws.Range('E42').Value = '  +19.15%  '
$ws.Range('D43').Value = "'12.06"
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = "'61.85"
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'9.08"
$ws.Range('E45').Value = '  +7.26%  '
$ws.Range('D46').Value = "'5.37"
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('E47').Value = '  +5.15%  '
$ws.Range('D48').Value = "'100.71"
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').Value = "'2.551.77"
$ws.Range('E50').Value = '  +4.26%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = "'0.182"
$ws.Range('E51').Value = '  +14.17%  '
